$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $find"
    }
    return $ok
}

# 1. Remove the second paragraph entirely ("This vignette is a guide to running ...").
#    It becomes merged away so the title paragraph is directly followed by the
#    (already empty) paragraph that used to be third.
$p2 = $d.Paragraphs(2)
$delRange = $d.Range($p2.Range.Start, $p2.Range.End)
$delRange.Delete()

# 2. Rewrite the title text.
Replace-Text 'Use of the “LabeledVolcanoPlot” DSP-DA plugin' 'Use of the “LabeledVolcanoPlot” DSP DA script'

# 3. Intended use paragraph rewrite.
Replace-Text 'This plug-in was designed for data from the GeoMx high-plex RNA assays, such as the CTA and WTA, or protein assays. It creates publication ready labeled volcano plots based on user inputs and statistical test results.' 'The LabeledVolcanoPlot DSP DA script was designed for data from the GeoMx nCounter (protein or RNA) or GeoMx NGS (CTA) readout applications. It creates publication ready labeled volcano plots based on user inputs and statistical test results.'

# 4. TOC + heading "Loading into the DSP-DA" -> "Loading into the DSP DA" (both occurrences).
Replace-Text 'Loading into the DSP-DA' 'Loading into the DSP DA'

# 5. "The LabeledVolcanoPlot plugin requires an extra file input from DSP-DA."
Replace-Text 'The LabeledVolcanoPlot plugin requires an extra file input from DSP-DA.' 'The LabeledVolcanoPlot script requires an extra file input from DSP DA.'

# 6. "creating a volcano plot in DSP-DA, the results file"
Replace-Text 'creating a volcano plot in DSP-DA, the results file' 'creating a volcano plot in DSP DA, the results file'

# 7. "before running the plugin."
Replace-Text 'before running the plugin.' 'before running the script.'

# 8. custom scripts section of the DSP-DA
Replace-Text '.R file may be loaded into the custom scripts section of the DSP-DA after you have a dataset processed and ready for analysis.' '.R file may be loaded into the custom scripts section of the DSP DA after you have a dataset processed and ready for analysis.'

# 9. de_results_filename description
Replace-Text 'Name of tab delimited file you’ve uploaded to the DSP-DA.' 'Name of tab delimited file you’ve uploaded to the DSP DA.'

# 10. Labels from DSP-DA volcano plot
Replace-Text 'Labels from DSP-DA volcano plot are not transferred to results file so must be user added' 'Labels from DSP DA volcano plot are not transferred to results file so must be user added'

# 11/12. Matching negative/positive axis labels
Replace-Text 'Matching negative (left) x-axis label to the volcano plot in DSP-DA' 'Matching negative (left) x-axis label to the volcano plot in DSP DA'
Replace-Text 'Matching positive (right) x-axis label to the volcano plot in DSP-DA' 'Matching positive (right) x-axis label to the volcano plot in DSP DA'

# 13. Output sentence in Example Parameter Set-up.
Replace-Text 'The LabeledVolcanoPlot plugin outputs a typical volcano plot figure with' 'The LabeledVolcanoPlot script outputs a typical volcano plot figure with'

Write-Output "Done"
